$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new SKU rows that were pasted in below the existing list ---
$ws.Range("A85").Value = 10004437
$ws.Range("A86").Value = 10126309
$ws.Range("A87").Value = 10025976

# --- A81 had picked up a stray "Open Sans" font/style (s="2") from a paste;
#     clear it back to the sheet's normal/default look used by every other cell ---
$ws.Range("A81").ClearFormats()

# --- Restore the on-screen selection to the new bottom of the list ---
$ws.Range("A81").Select()
